$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Prepend "OPENVAS: " to every Solution (column R) cell in the data rows (2-55).
for ($i = 2; $i -le 55; $i++) {
    $cell = $ws.Cells.Item($i, 18)
    $cell.Value2 = "OPENVAS: " + $cell.Value2
}

# 2. Add the new "scanner" header in column Z, matching the style of the other
#    header cells (bold, centered, bordered) by copying the format from Y1.
$ws.Range("Y1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1").Value = "scanner"

# 3. Touch the rest of column Z (rows 2-55) so the new column is materialised
#    across every data row, matching the expanded used range (A1:Z55).
$ws.Range("Z2:Z55").Style = "Normal"
